# Weekly data refresh: a new week's price record for Ciboulette at Vega
# Modelo de Temuco is inserted at row 153, pushing the existing rows
# 153-214 down to 154-215 (the sheet grows from 214 to 215 data rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 153, shifting everything
# below it (including formats) down by one row.
$ws.Rows(153).Insert()

# Populate the newly inserted row 153 with the new week's record.
$ws.Cells.Item(153, 1).Value = 10
$ws.Cells.Item(153, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(153, 3).Value = "La Araucanía"
$ws.Cells.Item(153, 4).Value = 44609
$ws.Cells.Item(153, 5).Value = 9
$ws.Cells.Item(153, 6).Value = 100112039
$ws.Cells.Item(153, 7).Value = "Ciboulette"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 40
$ws.Cells.Item(153, 11).Value = 5000
$ws.Cells.Item(153, 12).Value = 5000
$ws.Cells.Item(153, 13).Value = 5000
$ws.Cells.Item(153, 14).Value = "$/docena de atados"
$ws.Cells.Item(153, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(153, 16).Value = 1667
$ws.Cells.Item(153, 17).Value = 3
$ws.Cells.Item(153, 18).Value = "Hortaliza"
